# Helper function for reading all the split files.
# Adds two new sheets (splitFiles, joinFiles) modeled on the existing
# archiveFolder sheet, and repurposes the old archiveFolder sample paths.

$wb = $excel.ActiveWorkbook
$archiveSheet = $wb.Worksheets.Item("archiveFolder")

# --- Update archiveFolder's own sample data -------------------------------
# The old "/home/utkarsh/Desktop/exSitu" example is replaced with a new
# example path, and the two extra example rows (StarTracker / Test/png_avg)
# are cleared out since they are no longer needed as sample data.
$archiveSheet.Range("A2").Value = "/home/utkarsh/Desktop/restorePillars/LC2/01_HM_badChip_wet_bulk_converted"
$archiveSheet.Range("A3").Value = ""
$archiveSheet.Range("A4").Value = ""
$archiveSheet.Range("B4").Value = ""
$archiveSheet.Rows.Item(2).RowHeight = 25.5

# --- New sheet: splitFiles -------------------------------------------------
$splitSheet = $archiveSheet.Copy($null, $archiveSheet)
$splitSheet = $wb.Worksheets.Item($archiveSheet.Index + 1)
$splitSheet.Name = "splitFiles"

$splitSheet.Range("A1").Value = "File with path which you want to split"
$splitSheet.Range("B1").Value = "Delete flag (Set to 1 if you want to delete the directory after archiving. 0 otherwise. Default is 0)"
$splitSheet.Range("A2").Value = "/home/utkarsh/Projects/Datasets/Vehicles/Seattle Traffic in 5K 360° VR Video - Seattle Highways & Stadiums.webm"
$splitSheet.Range("B2").Value = 0
$splitSheet.Range("A3").Value = ""
$splitSheet.Range("A4").Value = ""
$splitSheet.Range("B4").Value = ""
$splitSheet.Rows.Item(2).RowHeight = 25.5

# --- New sheet: joinFiles ---------------------------------------------------
$joinSheet = $splitSheet.Copy($null, $splitSheet)
$joinSheet = $wb.Worksheets.Item($splitSheet.Index + 1)
$joinSheet.Name = "joinFiles"

$joinSheet.Range("A1").Value = "File with path which you want to join"
$joinSheet.Range("B1").Value = "Delete flag (Set to 1 if you want to delete the directory after archiving. 0 otherwise. Default is 0)"
$joinSheet.Range("A2").Value = "/home/utkarsh/Desktop/restorePillars/LC2/01_HM_badChip_wet_bulk_converted.zip_split_0001"
$joinSheet.Range("B2").Value = 1
$joinSheet.Rows.Item(2).RowHeight = 25.5

# --- Selection / active-tab bookkeeping ------------------------------------
$archiveSheet.Activate()
$archiveSheet.Range("B2").Select()

$splitSheet.Activate()
$splitSheet.Range("B3").Select()

$joinSheet.Activate()
$joinSheet.Range("B2").Select()

$splitSheet.Activate()

Write-Output "done"
